$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.99
$ws.Range("L2").Value = 1.53
$ws.Range("X2").Value = 10
$ws.Range("Y2").Value = 17
$ws.Range("Z2").Value = 50
$ws.Range("AC2").Value = 9.4
$ws.Range("AD2").Value = 28
$ws.Range("AF2").Value = 12.5
$ws.Range("AG2").Value = 13.5
$ws.Range("AH2").Value = 32
$ws.Range("AJ2").Value = 28
$ws.Range("AK2").Value = 32
$ws.Range("AL2").Value = 70
$ws.Range("G3").Value = 3.35
$ws.Range("J3").Value = 3.25
$ws.Range("P3").Value = 1.78
$ws.Range("Q3").Value = 2.02
$ws.Range("S3").Value = 3.6
$ws.Range("W3").Value = 1.42
$ws.Range("G4").Value = 2.06
$ws.Range("J4").Value = 2.88
$ws.Range("L4").Value = 1.54
$ws.Range("N4").Value = 2.22
$ws.Range("W4").Value = 1.94
$ws.Range("X4").Value = 9
$ws.Range("Y4").Value = 15.5
$ws.Range("Z4").Value = 55
$ws.Range("AB4").Value = 6.6
$ws.Range("AC4").Value = 9
$ws.Range("AD4").Value = 30
$ws.Range("AF4").Value = 11.5
$ws.Range("AG4").Value = 13.5
$ws.Range("AH4").Value = 36
$ws.Range("AJ4").Value = 27
$ws.Range("AK4").Value = 34
$ws.Range("AL4").Value = 85
$ws.Range("F5").Value = 1.43
$ws.Range("G5").Value = 1.51
$ws.Range("H5").Value = 1.04
$ws.Range("I5").Value = 14.5
$ws.Range("J5").Value = 3.7
$ws.Range("K5").Value = 4.6
$ws.Range("L5").Value = 1.45
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 2.54
$ws.Range("O5").Value = 1.42
$ws.Range("P5").Value = 1.57
$ws.Range("Q5").Value = 2.24
$ws.Range("R5").Value = 1.21
$ws.Range("S5").Value = 1.05
$ws.Range("T5").Value = 2.66
$ws.Range("U5").Value = 1.51
$ws.Range("V5").Value = 1.07
$ws.Range("W5").Value = 2.96
$ws.Range("X5").Value = 12
$ws.Range("Y5").Value = 28
$ws.Range("AB5").Value = 6.2
$ws.Range("AC5").Value = 12
$ws.Range("AD5").Value = 55
$ws.Range("AF5").Value = 8
$ws.Range("AG5").Value = 13
$ws.Range("AH5").Value = 50
$ws.Range("AJ5").Value = 14.5
$ws.Range("AK5").Value = 25
$ws.Range("AL5").Value = 85
$ws.Range("AN5").Value = 14
$ws.Range("G6").Value = 2.66
$ws.Range("H6").Value = 3.35
$ws.Range("K6").Value = 3.3
$ws.Range("L6").Value = 1.57
$ws.Range("W6").Value = 1.6
$ws.Range("Z6").Value = 28
$ws.Range("F7").Value = 1.51
$ws.Range("K7").Value = 8
$ws.Range("Q7").Value = 1.52
$ws.Range("R7").Value = 1.46
$ws.Range("T7").Value = 1.05
$ws.Range("U7").Value = 1.05
$ws.Range("S8").Value = 2.46
$ws.Range("Y8").Value = 15.5
$ws.Range("Z8").Value = 16
$ws.Range("AO8").Value = 9.800000000000001
$ws.Range("AF9").Value = 18
$ws.Range("F10").Value = 5.7
$ws.Range("L10").Value = 1.44
$ws.Range("O10").Value = 1.39
$ws.Range("P10").Value = 1.75
$ws.Range("Q10").Value = 2.12
$ws.Range("X10").Value = 14.5
$ws.Range("Y10").Value = 8.199999999999999
$ws.Range("Z10").Value = 11
$ws.Range("AA10").Value = 21
$ws.Range("AB10").Value = 21
$ws.Range("AC10").Value = 10
$ws.Range("AH10").Value = 28
$ws.Range("AI10").Value = 55
$ws.Range("K11").Value = 3.65
$ws.Range("F12").Value = 1.95
$ws.Range("G12").Value = 2.18
$ws.Range("H12").Value = 3.95
$ws.Range("I12").Value = 5.2
$ws.Range("J12").Value = 3.15
$ws.Range("K12").Value = 3.95
$ws.Range("L12").Value = 1.33
$ws.Range("M12").Value = 1.07
$ws.Range("N12").Value = 3.65
$ws.Range("O12").Value = 1.31
$ws.Range("P12").Value = 1.92
$ws.Range("Q12").Value = 1.76
$ws.Range("R12").Value = 1.35
$ws.Range("S12").Value = 3.3
$ws.Range("T12").Value = 1.78
$ws.Range("U12").Value = 2.08
$ws.Range("V12").Value = 1.28
$ws.Range("W12").Value = 1.85
$ws.Range("X12").Value = 18
$ws.Range("Y12").Value = 18.5
$ws.Range("AB12").Value = 11.5
$ws.Range("AC12").Value = 10
$ws.Range("AD12").Value = 21
$ws.Range("AF12").Value = 15.5
$ws.Range("AG12").Value = 13
$ws.Range("AN12").Value = 18
$ws.Range("F13").Value = 1.83
$ws.Range("L13").Value = 1.27
$ws.Range("P13").Value = 2.4
$ws.Range("S13").Value = 2.36
$ws.Range("X13").Value = 27
$ws.Range("AA13").Value = 95
$ws.Range("AB13").Value = 14.5
$ws.Range("AG13").Value = 12.5
$ws.Range("F14").Value = 2.6
$ws.Range("G14").Value = 2.86
$ws.Range("H14").Value = 2.72
$ws.Range("I14").Value = 3
$ws.Range("K14").Value = 3.75
$ws.Range("N14").Value = 3.3
$ws.Range("Q14").Value = 1.92
$ws.Range("S14").Value = 3.35
$ws.Range("T14").Value = 1.72
$ws.Range("U14").Value = 2.18
$ws.Range("W14").Value = 1.54
$ws.Range("I15").Value = 7.6
$ws.Range("J15").Value = 3.6
$ws.Range("L15").Value = 1.37
$ws.Range("S15").Value = 4
$ws.Range("X15").Value = 12.5
$ws.Range("G16").Value = 2.4
$ws.Range("H16").Value = 3.2
$ws.Range("I16").Value = 3.85
$ws.Range("J16").Value = 3.2
$ws.Range("K16").Value = 3.95
$ws.Range("L16").Value = 1.29
$ws.Range("N16").Value = 4.1
$ws.Range("P16").Value = 2.1
$ws.Range("Q16").Value = 1.75
$ws.Range("W16").Value = 1.72
$ws.Range("Z16").Value = 32
$ws.Range("AA16").Value = 75
$ws.Range("AE16").Value = 46
$ws.Range("AG16").Value = 14
$ws.Range("AI16").Value = 55
$ws.Range("AJ16").Value = 36
$ws.Range("AK16").Value = 28
$ws.Range("AO16").Value = 38
$ws.Range("G17").Value = 1.9
$ws.Range("S17").Value = 3.15
$ws.Range("W17").Value = 2.1
$ws.Range("AO17").Value = 65
$ws.Range("H18").Value = 2.24
$ws.Range("P18").Value = 1.87
$ws.Range("Q18").Value = 2.12
$ws.Range("AI18").Value = 40
$ws.Range("AK18").Value = 44
$ws.Range("AL18").Value = 60
$ws.Range("F19").Value = 3.15
$ws.Range("G19").Value = 4.9
$ws.Range("H19").Value = 1.97
$ws.Range("I19").Value = 2.44
$ws.Range("J19").Value = 3.15
$ws.Range("N19").Value = 2.9
$ws.Range("R19").Value = 1.25
$ws.Range("S19").Value = 2.26
$ws.Range("T19").Value = 1.05
$ws.Range("U19").Value = 1.05
$ws.Range("V19").Value = 1.7
$ws.Range("W19").Value = 1.26
$ws.Range("H20").Value = 3.4
$ws.Range("I20").Value = 5.6
$ws.Range("K20").Value = 4.9
$ws.Range("M20").Value = 1.09
$ws.Range("N20").Value = 2.7
$ws.Range("O20").Value = 1.41
$ws.Range("P20").Value = 1.66
$ws.Range("Q20").Value = 2.12
$ws.Range("R20").Value = 1.23
$ws.Range("S20").Value = 1.05
$ws.Range("T20").Value = 1.93
$ws.Range("U20").Value = 1.88
$ws.Range("Y20").Value = 18.5
$ws.Range("Z20").Value = 46
$ws.Range("AB20").Value = 11
$ws.Range("AC20").Value = 11
$ws.Range("AD20").Value = 26
$ws.Range("AF20").Value = 18
$ws.Range("AG20").Value = 16
$ws.Range("AH20").Value = 30
$ws.Range("AJ20").Value = 40
$ws.Range("AK20").Value = 38
$ws.Range("G21").Value = 5.5
$ws.Range("H21").Value = 1.82
$ws.Range("K21").Value = 3.85
$ws.Range("T21").Value = 1.86
$ws.Range("U21").Value = 1.96
$ws.Range("W21").Value = 1.22
$ws.Range("X21").Value = 15.5
$ws.Range("Y21").Value = 9.199999999999999
$ws.Range("AC21").Value = 8.800000000000001
$ws.Range("AF21").Value = 46
$ws.Range("AH21").Value = 24
$ws.Range("AL21").Value = 90
$ws.Range("G22").Value = 2.02
$ws.Range("H22").Value = 5.2
$ws.Range("J22").Value = 3.15
$ws.Range("U22").Value = 1.61
$ws.Range("V22").Value = 1.21
$ws.Range("W22").Value = 1.98
$ws.Range("X22").Value = 8.4
$ws.Range("AB22").Value = 7
$ws.Range("AC22").Value = 9.4
$ws.Range("AG22").Value = 14
$ws.Range("AH22").Value = 34
$ws.Range("AK22").Value = 34
$ws.Range("AL22").Value = 90
$ws.Range("AN22").Value = 34
